$wb = $excel.ActiveWorkbook

# --- branch sheet: fix r (resistance) values that were mistakenly left as
#     text "0" -> should be a tiny numeric value (nlp/data-entry fix) ---
$branch = $wb.Worksheets.Item("branch")

$branch.Range("E2").Value = 0.0000001
$branch.Range("E5").Value = 0.0000001
$branch.Range("E8").Value = 0.0000001

# --- add new "poles" column (N) for the MTDC branches ---
$header = $branch.Range("N1")
$header.Value = "poles"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.Item(7).LineStyle = 1
$header.Borders.Item(10).LineStyle = 1

$branch.Range("N2").Value = 2
$branch.Range("N3").Value = 2
$branch.Range("N4").Value = 2
$branch.Range("N5").Value = 2
$branch.Range("N6").Value = 2
$branch.Range("N7").Value = 2
$branch.Range("N8").Value = 2
$branch.Range("N9").Value = 2
$branch.Range("N10").Value = 2

# --- selection / active-sheet bookkeeping ---
$acLinks = $wb.Worksheets.Item("ac_links")
$acLinks.Range("O14").Select() | Out-Null

$branch.Activate()
$branch.Range("G13").Select() | Out-Null
